$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.285353422164917
$ws.Range("B1").Value = 2.410790205001831
$ws.Range("C1").Value = 3.237055063247681
$ws.Range("D1").Value = 3.437756776809692
$ws.Range("E1").Value = 1.06400203704834
